# Update profit/price figures across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Clear Glass Lens
$ws.Range("H33").Value = 22728498
$ws.Range("I33").Value = 31251324
$ws.Range("K33").Value = 31251324
$ws.Range("M33").Value = -31251095
# Row 51: Shark Oil
$ws.Range("H51").Value = 11564.294
$ws.Range("J51").Value = 7539.6
$ws.Range("L51").Value = 7539.6
$ws.Range("N51").Value = -8507.6
# Row 74: Wing Glue
$ws.Range("H74").Value = 88243390
$ws.Range("I74").Value = 136366700
$ws.Range("K74").Value = 136366700
$ws.Range("M74").Value = -136365764
# Row 77: Wing Glue
$ws.Range("H77").Value = 88243390
$ws.Range("I77").Value = 136366700
$ws.Range("K77").Value = 681833500
$ws.Range("M77").Value = -681828820
# Row 106: Enchanted Palladium Ink
$ws.Range("H106").Value = 200003800
$ws.Range("I106").Value = 333334660
$ws.Range("J106").Value = 7503
$ws.Range("K106").Value = 333334660
$ws.Range("L106").Value = 7503
$ws.Range("M106").Value = -333334029
$ws.Range("N106").Value = -8765
# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 1329.7646
$ws.Range("I132").Value = 1183.7084
$ws.Range("K132").Value = 3551.1252
$ws.Range("M132").Value = -1021.1252
# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 3160.6487
$ws.Range("J137").Value = 3323
$ws.Range("L137").Value = 9969
$ws.Range("N137").Value = -15069
# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2047469.1
$ws.Range("J138").Value = 2507877.5
$ws.Range("L138").Value = 7523632.5
$ws.Range("N138").Value = -7533912.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 1426681.9
$ws.Range("I32").Value = 1529424.2
$ws.Range("K32").Value = 1529424.2
$ws.Range("M32").Value = -1529137.2

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Iron Ingot
$ws.Range("H20").Value = 9806157
$ws.Range("I20").Value = 16668992
$ws.Range("J20").Value = 2107.1428
$ws.Range("K20").Value = 16668992
$ws.Range("L20").Value = 2107.1428
$ws.Range("M20").Value = -16668745
$ws.Range("N20").Value = -2601.1428
# Row 107: Deepgold Nugget
$ws.Range("H107").Value = 45002956
$ws.Range("I107").Value = 62501332
$ws.Range("K107").Value = 62501332
$ws.Range("M107").Value = -62499412
# Row 109: Deepgold Sledgehammer
$ws.Range("H109").Value = 55938
$ws.Range("J109").Value = 55938
$ws.Range("L109").Value = 55938
$ws.Range("N109").Value = -58712

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Ash Lumber
$ws.Range("H16").Value = 5412.5
$ws.Range("I16").Value = 3143.5
$ws.Range("J16").Value = 6773.9
$ws.Range("K16").Value = 3143.5
$ws.Range("L16").Value = 6773.9
$ws.Range("M16").Value = -2856.5
$ws.Range("N16").Value = -7347.9
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 7014.654
$ws.Range("I31").Value = 3228.875
$ws.Range("K31").Value = 3228.875
$ws.Range("M31").Value = -2933.875
# Row 34: Walnut Lumber
$ws.Range("H34").Value = 7014.654
$ws.Range("I34").Value = 3228.875
$ws.Range("K34").Value = 3228.875
$ws.Range("M34").Value = -3026.875
# Row 58: Mahogany Lumber
$ws.Range("H58").Value = 14293093
$ws.Range("I58").Value = 33335988
$ws.Range("K58").Value = 33335988
$ws.Range("M58").Value = -33335785
# Row 94: Beech Lumber
$ws.Range("H94").Value = 819.4091
$ws.Range("I94").Value = 981
$ws.Range("J94").Value = 771.8823
$ws.Range("K94").Value = 981
$ws.Range("L94").Value = 771.8823
$ws.Range("M94").Value = -530
$ws.Range("N94").Value = -1673.8823
# Row 113: White Ash Lumber
$ws.Range("H113").Value = 5412.5
$ws.Range("I113").Value = 3143.5
$ws.Range("J113").Value = 6773.9
$ws.Range("K113").Value = 3143.5
$ws.Range("L113").Value = 6773.9
$ws.Range("M113").Value = -973.5
$ws.Range("N113").Value = -11113.9
# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 7040.1816
$ws.Range("I132").Value = 3388.625
$ws.Range("J132").Value = 10476.941
$ws.Range("K132").Value = 10165.875
$ws.Range("L132").Value = 31430.823
$ws.Range("M132").Value = -7635.875
$ws.Range("N132").Value = -36490.823
# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 5118.25
$ws.Range("I134").Value = 2243.5
$ws.Range("K134").Value = 6730.5
$ws.Range("M134").Value = -4195.5
# Row 136: Dark Mahogany Lumber
$ws.Range("H136").Value = 14293093
$ws.Range("I136").Value = 33335988
$ws.Range("K136").Value = 100007964
$ws.Range("M136").Value = -100005414

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Table Salt
$ws.Range("H2").Value = 69369.44500000001
$ws.Range("I2").Value = 50.142857
$ws.Range("J2").Value = 251332.62
$ws.Range("K2").Value = 300.857142
$ws.Range("L2").Value = 1507995.72
$ws.Range("M2").Value = -187.857142
$ws.Range("N2").Value = -1508221.72
# Row 107: Frantoio Oil
$ws.Range("H107").Value = 40000420
$ws.Range("I107").Value = 599
$ws.Range("J107").Value = 50000376
$ws.Range("K107").Value = 1797
$ws.Range("L107").Value = 150001128
$ws.Range("M107").Value = 123
$ws.Range("N107").Value = -150004968

$ws = $wb.Worksheets.Item("GSM")
# Row 31: Staghorn Staff
$ws.Range("H31").Value = 2218
$ws.Range("I31").Value = 191
$ws.Range("J31").Value = 2724.75
$ws.Range("K31").Value = 191
$ws.Range("L31").Value = 2724.75
$ws.Range("M31").Value = 101
$ws.Range("N31").Value = -3308.75
# Row 37: Toothed Staghorn Staff
$ws.Range("H37").Value = 2218
$ws.Range("I37").Value = 191
$ws.Range("J37").Value = 2724.75
$ws.Range("K37").Value = 191
$ws.Range("L37").Value = 2724.75
$ws.Range("M37").Value = 86
$ws.Range("N37").Value = -3278.75

$ws = $wb.Worksheets.Item("LTW")
# Row 14: Hard Leather Shoes
$ws.Range("H14").Value = 30000
$ws.Range("I14").Value = 30000
$ws.Range("K14").Value = 30000
$ws.Range("M14").Value = -29828
# Row 55: Peiste Leather
$ws.Range("H55").Value = 349.81482
$ws.Range("J55").Value = 542.4666999999999
$ws.Range("L55").Value = 542.4666999999999
$ws.Range("N55").Value = -888.4666999999999
# Row 68: Wyvern Leather
$ws.Range("H68").Value = 4676.375
$ws.Range("I68").Value = 3324.8333
$ws.Range("K68").Value = 3324.8333
$ws.Range("M68").Value = -2575.8333
# Row 71: Wyvern Leather
$ws.Range("H71").Value = 4676.375
$ws.Range("I71").Value = 3324.8333
$ws.Range("K71").Value = 16624.1665
$ws.Range("M71").Value = -12880.1665
# Row 100: Tiger Leather
$ws.Range("H100").Value = 2843.5
$ws.Range("I100").Value = 1873.5
$ws.Range("K100").Value = 1873.5
$ws.Range("M100").Value = -1332.5
# Row 136: Br'aax Leather
$ws.Range("H136").Value = 10260
$ws.Range("I136").Value = 2280
$ws.Range("K136").Value = 6840
$ws.Range("M136").Value = -4290

$ws = $wb.Worksheets.Item("WVR")
# Row 8: Hempen Tabard
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
# Row 10: Hempen Kecks
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

